$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Force the cell to remain text (matches the source inlineStr cells)
    # even when the literal looks like a pure number (e.g. "5.48").
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "64.698.73"
Set-TextValue $ws.Range("E2") "  +0.41%  "

Set-TextValue $ws.Range("D3") "3.160.87"
Set-TextValue $ws.Range("E3") "  +0.44%  "

Set-TextValue $ws.Range("E4") "  +0.19%  "

Set-TextValue $ws.Range("D5") "617.22"
Set-TextValue $ws.Range("E5") "  +2.30%  "

Set-TextValue $ws.Range("D6") "146.96"
Set-TextValue $ws.Range("E6") "  -2.60%  "

Set-TextValue $ws.Range("E7") "  +0.06%  "

Set-TextValue $ws.Range("D8") "3.154.97"
Set-TextValue $ws.Range("E8") "  +0.25%  "

Set-TextValue $ws.Range("E9") "  -1.25%  "

Set-TextValue $ws.Range("E10") "  -1.37%  "

Set-TextValue $ws.Range("D11") "5.48"
Set-TextValue $ws.Range("E11") "  -2.52%  "

Set-TextValue $ws.Range("E12") "  -1.34%  "

Set-TextValue $ws.Range("D13") "0.0000259"
Set-TextValue $ws.Range("E13") "  -0.93%  "

Set-TextValue $ws.Range("D14") "35.77"
Set-TextValue $ws.Range("E14") "  -3.62%  "

Set-TextValue $ws.Range("D15") "3.684.68"
Set-TextValue $ws.Range("E15") "  +1.76%  "

Set-TextValue $ws.Range("E16") "  +2.71%  "

Set-TextValue $ws.Range("D17") "64.679.41"
Set-TextValue $ws.Range("E17") "  +0.27%  "

Set-TextValue $ws.Range("D18") "3.165.50"
Set-TextValue $ws.Range("E18") "  +0.94%  "

Set-TextValue $ws.Range("D19") "6.92"
Set-TextValue $ws.Range("E19") "  -2.10%  "

Set-TextValue $ws.Range("D20") "478.49"
Set-TextValue $ws.Range("E20") "  -1.21%  "

Set-TextValue $ws.Range("D21") "14.71"
Set-TextValue $ws.Range("E21") "  -0.56%  "

Set-TextValue $ws.Range("E22") "  +1.04%  "

Set-TextValue $ws.Range("D23") "7.93"
Set-TextValue $ws.Range("E23") "  +1.16%  "

Set-TextValue $ws.Range("D24") "13.78"
Set-TextValue $ws.Range("E24") "  -1.45%  "

Set-TextValue $ws.Range("D25") "84.34"
Set-TextValue $ws.Range("E25") "  -0.57%  "

Set-TextValue $ws.Range("E26") "  -0.02%  "

Set-TextValue $ws.Range("E27") "  -3.68%  "

Set-TextValue $ws.Range("D28") "8.56"
Set-TextValue $ws.Range("E28") "  -2.08%  "

Set-TextValue $ws.Range("B29") "Hedera"
Set-TextValue $ws.Range("C29") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D29") "0.118"
Set-TextValue $ws.Range("E29") "  -6.47%  "

Set-TextValue $ws.Range("B30") "NEARProtocol"
Set-TextValue $ws.Range("C30") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D30") "6.91"
Set-TextValue $ws.Range("E30") "  -3.83%  "

Set-TextValue $ws.Range("E31") "  -8.49%  "

Set-TextValue $ws.Range("E32") "  +0.30%  "

Set-TextValue $ws.Range("E33") "  -1.32%  "

Set-TextValue $ws.Range("D34") "26.54"
Set-TextValue $ws.Range("E34") "  -1.35%  "

Set-TextValue $ws.Range("E35") "  +1.87%  "

Set-TextValue $ws.Range("D36") "0.0₃0778"
Set-TextValue $ws.Range("E36") "  +1.86%  "

Set-TextValue $ws.Range("E37") "  -2.14%  "

Set-TextValue $ws.Range("D38") "53.05"
Set-TextValue $ws.Range("E38") "  -3.00%  "

Set-TextValue $ws.Range("D39") "3.16"
Set-TextValue $ws.Range("E39") "  -5.09%  "

Set-TextValue $ws.Range("D40") "460.14"
Set-TextValue $ws.Range("E40") "  +1.40%  "

Set-TextValue $ws.Range("E41") "  -1.41%  "

Set-TextValue $ws.Range("E42") "  -4.18%  "

Set-TextValue $ws.Range("D43") "8.39"
Set-TextValue $ws.Range("E43") "  -1.98%  "

Set-TextValue $ws.Range("D44") "2.846.71"
Set-TextValue $ws.Range("E44") "  -2.06%  "

Set-TextValue $ws.Range("E45") "  -4.92%  "

Set-TextValue $ws.Range("D46") "0.268"
Set-TextValue $ws.Range("E46") "  -3.00%  "

Set-TextValue $ws.Range("D47") "2.42"
Set-TextValue $ws.Range("E47") "  +3.12%  "

Set-TextValue $ws.Range("D48") "26.55"
Set-TextValue $ws.Range("E48") "  -1.94%  "

Set-TextValue $ws.Range("E49") "  +0.14%  "

Set-TextValue $ws.Range("E50") "  -2.04%  "

Set-TextValue $ws.Range("D51") "120.34"
Set-TextValue $ws.Range("E51") "  +0.18%  "
